# Auto-generated update of Leve profit/price figures (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 626.8570999999999  # ALC!H6
$ws.Cells.Item(6, 9).Value = 626.8570999999999  # ALC!I6
$ws.Cells.Item(6, 10).Value = 0  # ALC!J6
$ws.Cells.Item(6, 11).Value = 1880.5713  # ALC!K6
$ws.Cells.Item(6, 12).Value = 0  # ALC!L6
$ws.Cells.Item(6, 13).ClearContents()  # ALC!M6
$ws.Cells.Item(6, 14).Value = -1768.5713  # ALC!N6

$ws.Cells.Item(9, 8).Value = 236.5  # ALC!H9
$ws.Cells.Item(9, 9).Value = 243.8  # ALC!I9
$ws.Cells.Item(9, 10).Value = 200  # ALC!J9
$ws.Cells.Item(9, 11).Value = 243.8  # ALC!K9
$ws.Cells.Item(9, 12).Value = 200  # ALC!L9
$ws.Cells.Item(9, 13).Value = -74.80000000000001  # ALC!M9
$ws.Cells.Item(9, 14).Value = -538  # ALC!N9

$ws.Cells.Item(12, 8).Value = 934.7059  # ALC!H12
$ws.Cells.Item(12, 9).Value = 934.7059  # ALC!I12
$ws.Cells.Item(12, 10).Value = 0  # ALC!J12
$ws.Cells.Item(12, 11).Value = 934.7059  # ALC!K12
$ws.Cells.Item(12, 12).Value = 0  # ALC!L12
$ws.Cells.Item(12, 13).ClearContents()  # ALC!M12
$ws.Cells.Item(12, 14).Value = -764.7059  # ALC!N12

$ws.Cells.Item(132, 8).Value = 1579.4  # ALC!H132
$ws.Cells.Item(132, 9).Value = 1253.9056  # ALC!I132
$ws.Cells.Item(132, 10).Value = 10205  # ALC!J132
$ws.Cells.Item(132, 11).Value = 3761.7168  # ALC!K132
$ws.Cells.Item(132, 12).Value = 30615  # ALC!L132
$ws.Cells.Item(132, 13).Value = -1231.7168  # ALC!M132
$ws.Cells.Item(132, 14).Value = -35675  # ALC!N132

$ws.Cells.Item(137, 8).Value = 4404.1724  # ALC!H137
$ws.Cells.Item(137, 9).Value = 3603.2354  # ALC!I137
$ws.Cells.Item(137, 11).Value = 10809.7062  # ALC!K137
$ws.Cells.Item(137, 13).Value = -8259.706200000001  # ALC!M137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1352.4127  # ARM!H32
$ws.Cells.Item(32, 9).Value = 1417.2456  # ARM!I32
$ws.Cells.Item(32, 11).Value = 1417.2456  # ARM!K32
$ws.Cells.Item(32, 13).Value = -1130.2456  # ARM!M32

$ws.Cells.Item(61, 8).Value = 8052.2856  # ARM!H61
$ws.Cells.Item(61, 9).Value = 6451.1665  # ARM!I61
$ws.Cells.Item(61, 10).Value = 10187.111  # ARM!J61
$ws.Cells.Item(61, 11).Value = 6451.1665  # ARM!K61
$ws.Cells.Item(61, 12).Value = 10187.111  # ARM!L61
$ws.Cells.Item(61, 13).Value = -6239.1665  # ARM!M61
$ws.Cells.Item(61, 14).Value = -10611.111  # ARM!N61

$ws.Cells.Item(74, 8).Value = 6947992  # ARM!H74
$ws.Cells.Item(74, 9).Value = 8550113  # ARM!I74
$ws.Cells.Item(74, 10).Value = 5468.5557  # ARM!J74
$ws.Cells.Item(74, 11).Value = 8550113  # ARM!K74
$ws.Cells.Item(74, 12).Value = 5468.5557  # ARM!L74
$ws.Cells.Item(74, 13).Value = -8549239  # ARM!M74
$ws.Cells.Item(74, 14).Value = -7216.5557  # ARM!N74

$ws.Cells.Item(77, 8).Value = 6947992  # ARM!H77
$ws.Cells.Item(77, 9).Value = 8550113  # ARM!I77
$ws.Cells.Item(77, 10).Value = 5468.5557  # ARM!J77
$ws.Cells.Item(77, 11).Value = 42750565  # ARM!K77
$ws.Cells.Item(77, 12).Value = 27342.7785  # ARM!L77
$ws.Cells.Item(77, 13).Value = -42746197  # ARM!M77
$ws.Cells.Item(77, 14).Value = -36078.7785  # ARM!N77

$ws.Cells.Item(102, 8).Value = 1838.2858  # ARM!H102
$ws.Cells.Item(102, 9).Value = 1787.3846  # ARM!I102
$ws.Cells.Item(102, 11).Value = 1787.3846  # ARM!K102
$ws.Cells.Item(102, 13).Value = -165.3846000000001  # ARM!M102

$ws.Cells.Item(122, 8).Value = 3514.8538  # ARM!H122
$ws.Cells.Item(122, 9).Value = 3135  # ARM!I122
$ws.Cells.Item(122, 11).Value = 9405  # ARM!K122
$ws.Cells.Item(122, 13).Value = -6955  # ARM!M122

$ws.Cells.Item(132, 8).Value = 6176.316  # ARM!H132
$ws.Cells.Item(132, 9).Value = 4501.8623  # ARM!I132
$ws.Cells.Item(132, 10).Value = 11571.777  # ARM!J132
$ws.Cells.Item(132, 11).Value = 13505.5869  # ARM!K132
$ws.Cells.Item(132, 12).Value = 34715.331  # ARM!L132
$ws.Cells.Item(132, 13).Value = -10975.5869  # ARM!M132
$ws.Cells.Item(132, 14).Value = -39775.331  # ARM!N132

$ws.Cells.Item(136, 8).Value = 8052.2856  # ARM!H136
$ws.Cells.Item(136, 9).Value = 6451.1665  # ARM!I136
$ws.Cells.Item(136, 10).Value = 10187.111  # ARM!J136
$ws.Cells.Item(136, 11).Value = 19353.4995  # ARM!K136
$ws.Cells.Item(136, 12).Value = 30561.333  # ARM!L136
$ws.Cells.Item(136, 13).Value = -16803.4995  # ARM!M136
$ws.Cells.Item(136, 14).Value = -35661.333  # ARM!N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3165.25  # BSM!H86
$ws.Cells.Item(86, 9).Value = 1687.875  # BSM!I86
$ws.Cells.Item(86, 10).Value = 6120  # BSM!J86
$ws.Cells.Item(86, 11).Value = 1687.875  # BSM!K86
$ws.Cells.Item(86, 12).Value = 6120  # BSM!L86
$ws.Cells.Item(86, 13).Value = -564.875  # BSM!M86
$ws.Cells.Item(86, 14).Value = -8366  # BSM!N86

$ws.Cells.Item(89, 8).Value = 3165.25  # BSM!H89
$ws.Cells.Item(89, 9).Value = 1687.875  # BSM!I89
$ws.Cells.Item(89, 10).Value = 6120  # BSM!J89
$ws.Cells.Item(89, 11).Value = 8439.375  # BSM!K89
$ws.Cells.Item(89, 12).Value = 30600  # BSM!L89
$ws.Cells.Item(89, 13).Value = -2823.375  # BSM!M89
$ws.Cells.Item(89, 14).Value = -41832  # BSM!N89

$ws.Cells.Item(99, 8).Value = 2961  # BSM!H99
$ws.Cells.Item(99, 9).Value = 2959.4  # BSM!I99
$ws.Cells.Item(99, 10).Value = 2966.3333  # BSM!J99
$ws.Cells.Item(99, 11).Value = 2959.4  # BSM!K99
$ws.Cells.Item(99, 12).Value = 2966.3333  # BSM!L99
$ws.Cells.Item(99, 13).Value = -1461.4  # BSM!M99
$ws.Cells.Item(99, 14).Value = -5962.3333  # BSM!N99

$ws.Cells.Item(105, 8).Value = 18800.066  # BSM!H105
$ws.Cells.Item(105, 9).Value = 19502.666  # BSM!I105
$ws.Cells.Item(105, 11).Value = 19502.666  # BSM!K105
$ws.Cells.Item(105, 13).Value = -17755.666  # BSM!M105

$ws.Cells.Item(107, 8).Value = 4073.875  # BSM!H107
$ws.Cells.Item(107, 9).Value = 3848.5  # BSM!I107
$ws.Cells.Item(107, 11).Value = 3848.5  # BSM!K107
$ws.Cells.Item(107, 13).Value = -1928.5  # BSM!M107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1589.3572  # CRP!H22
$ws.Cells.Item(22, 9).Value = 895.6818  # CRP!I22
$ws.Cells.Item(22, 11).Value = 895.6818  # CRP!K22
$ws.Cells.Item(22, 13).Value = -545.6818  # CRP!M22

$ws.Cells.Item(58, 8).Value = 4839.56  # CRP!H58
$ws.Cells.Item(58, 9).Value = 2426.1428  # CRP!I58
$ws.Cells.Item(58, 10).Value = 7911.1816  # CRP!J58
$ws.Cells.Item(58, 11).Value = 2426.1428  # CRP!K58
$ws.Cells.Item(58, 12).Value = 7911.1816  # CRP!L58
$ws.Cells.Item(58, 13).Value = -2223.1428  # CRP!M58
$ws.Cells.Item(58, 14).Value = -8317.1816  # CRP!N58

$ws.Cells.Item(107, 8).Value = 2111.4119  # CRP!H107
$ws.Cells.Item(107, 9).Value = 1684.3077  # CRP!I107
$ws.Cells.Item(107, 11).Value = 1684.3077  # CRP!K107
$ws.Cells.Item(107, 13).Value = 235.6922999999999  # CRP!M107

$ws.Cells.Item(134, 8).Value = 2221.1191  # CRP!H134
$ws.Cells.Item(134, 9).Value = 1323.3823  # CRP!I134
$ws.Cells.Item(134, 10).Value = 6036.5  # CRP!J134
$ws.Cells.Item(134, 11).Value = 3970.1469  # CRP!K134
$ws.Cells.Item(134, 12).Value = 18109.5  # CRP!L134
$ws.Cells.Item(134, 13).Value = -1435.1469  # CRP!M134
$ws.Cells.Item(134, 14).Value = -23179.5  # CRP!N134

$ws.Cells.Item(136, 8).Value = 4839.56  # CRP!H136
$ws.Cells.Item(136, 9).Value = 2426.1428  # CRP!I136
$ws.Cells.Item(136, 10).Value = 7911.1816  # CRP!J136
$ws.Cells.Item(136, 11).Value = 7278.428400000001  # CRP!K136
$ws.Cells.Item(136, 12).Value = 23733.5448  # CRP!L136
$ws.Cells.Item(136, 13).Value = -4728.428400000001  # CRP!M136
$ws.Cells.Item(136, 14).Value = -28833.5448  # CRP!N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 66.95999999999999  # CUL!H33
$ws.Cells.Item(33, 9).Value = 63.42857  # CUL!I33
$ws.Cells.Item(33, 10).Value = 85.5  # CUL!J33
$ws.Cells.Item(33, 11).Value = 380.57142  # CUL!K33
$ws.Cells.Item(33, 12).Value = 513  # CUL!L33
$ws.Cells.Item(33, 13).Value = -97.57141999999999  # CUL!M33
$ws.Cells.Item(33, 14).Value = -1079  # CUL!N33

$ws.Cells.Item(107, 8).Value = 2150.375  # CUL!H107
$ws.Cells.Item(107, 9).Value = 399.8  # CUL!I107
$ws.Cells.Item(107, 10).Value = 5068  # CUL!J107
$ws.Cells.Item(107, 11).Value = 1199.4  # CUL!K107
$ws.Cells.Item(107, 12).Value = 15204  # CUL!L107
$ws.Cells.Item(107, 13).Value = 720.5999999999999  # CUL!M107
$ws.Cells.Item(107, 14).Value = -19044  # CUL!N107

$ws.Cells.Item(132, 8).Value = 3527.2856  # CUL!H132
$ws.Cells.Item(132, 9).Value = 2687.5  # CUL!I132
$ws.Cells.Item(132, 10).Value = 4290.727  # CUL!J132
$ws.Cells.Item(132, 11).Value = 24187.5  # CUL!K132
$ws.Cells.Item(132, 12).Value = 38616.543  # CUL!L132
$ws.Cells.Item(132, 13).Value = -21657.5  # CUL!M132
$ws.Cells.Item(132, 14).Value = -43676.543  # CUL!N132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(131, 8).Value = 44154  # GSM!H131
$ws.Cells.Item(131, 9).Value = 28987  # GSM!I131
$ws.Cells.Item(131, 10).Value = 47945.75  # GSM!J131
$ws.Cells.Item(131, 11).Value = 28987  # GSM!K131
$ws.Cells.Item(131, 12).Value = 47945.75  # GSM!L131
$ws.Cells.Item(131, 13).Value = -23947  # GSM!M131
$ws.Cells.Item(131, 14).Value = -58025.75  # GSM!N131

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4350.378  # LTW!H132
$ws.Cells.Item(132, 9).Value = 3027.3125  # LTW!I132
$ws.Cells.Item(132, 10).Value = 7607.154  # LTW!J132
$ws.Cells.Item(132, 11).Value = 9081.9375  # LTW!K132
$ws.Cells.Item(132, 12).Value = 22821.462  # LTW!L132
$ws.Cells.Item(132, 13).Value = -6551.9375  # LTW!M132
$ws.Cells.Item(132, 14).Value = -27881.462  # LTW!N132

$ws.Cells.Item(136, 8).Value = 7182.5  # LTW!H136
$ws.Cells.Item(136, 9).Value = 4793.2  # LTW!I136
$ws.Cells.Item(136, 11).Value = 14379.6  # LTW!K136
$ws.Cells.Item(136, 13).Value = -11829.6  # LTW!M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 9400  # WVR!H81
$ws.Cells.Item(81, 10).Value = 18652  # WVR!J81
$ws.Cells.Item(81, 12).Value = 37304  # WVR!L81
$ws.Cells.Item(81, 14).Value = -39426  # WVR!N81

$ws.Cells.Item(84, 8).Value = 9400  # WVR!H84
$ws.Cells.Item(84, 10).Value = 18652  # WVR!J84
$ws.Cells.Item(84, 12).Value = 186520  # WVR!L84
$ws.Cells.Item(84, 14).Value = -197128  # WVR!N84

$ws.Cells.Item(96, 8).Value = 3666.6667  # WVR!H96
$ws.Cells.Item(96, 10).Value = 4000  # WVR!J96
$ws.Cells.Item(96, 12).Value = 4000  # WVR!L96
$ws.Cells.Item(96, 14).Value = -6746  # WVR!N96

$ws.Cells.Item(100, 8).Value = 800.1111  # WVR!H100
$ws.Cells.Item(100, 9).Value = 605.6  # WVR!I100
$ws.Cells.Item(100, 11).Value = 1211.2  # WVR!K100
$ws.Cells.Item(100, 13).Value = -670.2  # WVR!M100

$ws.Cells.Item(122, 8).Value = 2555  # WVR!H122
$ws.Cells.Item(122, 9).Value = 1470.1666  # WVR!I122
$ws.Cells.Item(122, 10).Value = 13403.333  # WVR!J122
$ws.Cells.Item(122, 11).Value = 4410.4998  # WVR!K122
$ws.Cells.Item(122, 12).Value = 40209.999  # WVR!L122
$ws.Cells.Item(122, 13).Value = -1960.4998  # WVR!M122
$ws.Cells.Item(122, 14).Value = -45109.999  # WVR!N122
